$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple cell value edits (rows 1-25 keep their row numbers) ---

# Row 3: E3 -5.7 -> blank
$ws.Range("E3").ClearContents()

# Row 5: F5 17.66 -> blank
$ws.Range("F5").ClearContents()

# Row 8: F8 blank -> 17.05
$ws.Range("F8").Value = 17.05

# Row 10: F10 blank -> 16.43
$ws.Range("F10").Value = 16.43

# Row 12: F12 17.45 -> blank
$ws.Range("F12").ClearContents()

# Row 15: F15 blank -> 16.2
$ws.Range("F15").Value = 16.2

# Row 18: F18 18.35 -> blank
$ws.Range("F18").ClearContents()

# Row 19: F19 17.81 -> blank
$ws.Range("F19").ClearContents()

# Row 25: F25 blank -> 16.6
$ws.Range("F25").Value = 16.6

# --- Remove whole rows (delete in descending order so row numbers stay valid) ---
# Row 28 = "SC 92" is removed entirely
$ws.Rows(28).Delete()
# Row 26 = "RM 232" is removed entirely
$ws.Rows(26).Delete()

# --- After the deletions, rows 27-35 shifted up by two; fix remaining value diffs ---
# New row 26 ("SC 5"): D26 blank -> -13.8
$ws.Range("D26").Value = -13.8

# New row 27 ("SC 101"): D27 -14.6 -> blank
$ws.Range("D27").ClearContents()

# New row 29 ("SC 119"): F29 18.06 -> blank
$ws.Range("F29").ClearContents()

# New row 33 ("SC 232"): D33 blank -> -14.1, E33 blank -> -10.7
$ws.Range("D33").Value = -14.1
$ws.Range("E33").Value = -10.7
